# Battery.xlsx update: add Tenpower IFR26700-45HE / IFR26700-40HE rows,
# bump quantity-per-device from 10 to 25, highlight the new computed
# price cells in yellow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("26650")

# 1) Quantity per device: 10 -> 25 (this ripples through V/W/X columns)
$ws.Range("V2").Value = 25

# 2) New row 8: Tenpower IFR26700-45HE
$ws.Range("A8").Value = "NKON"
$ws.Range("B8").Value = "Tenpower"
$ws.Range("C8").Value = "IFR26700-45HE"
$ws.Range("D8").Value = 3.2
$ws.Range("E8").Value = 2.5
$ws.Range("F8").Value = 4.45
$ws.Range("G8").Value = 4.5
$ws.Range("H8").Value = 9
$ws.Range("M8").Value = 1.95
$ws.Range("N8").Value = 1.85

# 3) New row 9: Tenpower IFR26700-40HE
$ws.Range("A9").Value = "NKON"
$ws.Range("B9").Value = "Tenpower"
$ws.Range("C9").Value = "IFR26700-40HE"
$ws.Range("D9").Value = 3.2
$ws.Range("E9").Value = 2.5
$ws.Range("F9").Value = 3.95
$ws.Range("G9").Value = 4
$ws.Range("H9").Value = 8
$ws.Range("M9").Value = 1.85
$ws.Range("N9").Value = 1.75

# 4) Price-tier formulas for the new rows (Tenpower only quoted the 100pc
#    price, so the lower tiers are back-computed from it)
$ws.Range("I8").Formula = '=$M8*$M$2/I$2'
$ws.Range("I9").Formula = '=$M9*$M$2/I$2'
$ws.Range("J8:L9").Formula = '=$M8*$M$2/J$2'

# 5) Highlight the computed price cells in yellow
$ws.Range("I8:L9").Interior.Color = 65535

# 6) Extend the existing calc columns down through the new rows
$ws.Range("R8:R9").Formula = '=ROUNDUP(R$2/D8/H8,2)'
$ws.Range("T8:T9").Formula = '=ROUNDUP(T$2/$D8/$F8,2)'
$ws.Range("U8:U9").Formula = '=ROUNDUP(MAX(T$2/$D8/$F8,R$2/E8/H8),0)'
$ws.Range("V8:V9").Formula = '=V$2*$U8'
$ws.Range("W8:W9").Formula = '=IF(AND(V8>=P$2,ISNUMBER(P8)),V8*P8,IF(AND(V8>=O$2,ISNUMBER(O8)),V8*O8,IF(AND(V8>=N$2,ISNUMBER(N8)),V8*N8,IF(AND(V8>=M$2,ISNUMBER(M8)),V8*M8,IF(AND(V8>=L$2,ISNUMBER(L8)),V8*L8,IF(AND(V8>=K$2,ISNUMBER(K8)),V8*K8,IF(AND(V8>=J$2,ISNUMBER(J8)),V8*J8,V8*I8)))))))'
$ws.Range("X8:X9").Formula = '=W8/V8*U8'
$ws.Range("Y8:Y9").Formula = '=D8*F8*U8'
$ws.Range("Z8:Z9").Formula = '=Y8/R$2*60'
$ws.Range("AA8:AA9").Formula = '=E8*H8*U8'

# 7) Hyperlinks for the new type cells, matching the existing "Link" style
$ws.Hyperlinks.Add($ws.Range("C8"), "https://www.nkon.nl/de/rechargeable/lifepo4/26700/tenpower-ifr26700-45he-4500mah-9a-lifepo4.html") | Out-Null
$ws.Range("C8").Style = "Link"
$ws.Hyperlinks.Add($ws.Range("C9"), "https://www.nkon.nl/de/rechargeable/lifepo4/26700/tenpower-ifr26700-40he-4000mah-8a-lifepo4.html") | Out-Null
$ws.Range("C9").Style = "Link"

# 8) Update the current selection to reflect where the author ended up
$ws.Range("X5").Select() | Out-Null
